$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (old row 7) - the new table only has 5 data rows
$ws.Rows.Item(7).Delete()

# Insert a new blank column at the start; this shifts the existing
# Depth from / Depth to / Unit weight / Soil type columns from A:D to B:E
$ws.Columns.Item(1).Insert()

# Header row: new column A has no header label (leave A1 blank)

# New index column (A2:A6) = 0..4, formatted like the header row (bold, bordered, centered)
$ws.Range("B1").Copy()
$ws.Range("A2:A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4

# Replace the depth / unit weight / soil type data with the new dataset
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 3.16
$ws.Range("D2").Value = 18
$ws.Range("E2").Value = "SAND"

$ws.Range("B3").Value = 3.16
$ws.Range("C3").Value = 5.9
$ws.Range("D3").Value = 17
$ws.Range("E3").Value = "CLAY"

$ws.Range("B4").Value = 5.9
$ws.Range("C4").Value = 14.86
$ws.Range("D4").Value = 19.5
$ws.Range("E4").Value = "SAND"

$ws.Range("B5").Value = 14.86
$ws.Range("C5").Value = 15.7
$ws.Range("D5").Value = 20
$ws.Range("E5").Value = "SAND"

$ws.Range("B6").Value = 15.7
$ws.Range("C6").Value = 20
$ws.Range("D6").Value = 20
$ws.Range("E6").Value = "SAND"
